$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.024.79"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.502.95"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.20"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.01"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.497.07"
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.193"
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.27"
$ws.Range("E11").Value = "  +7.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.579"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.03"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000274"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.062.92"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "611.97"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.25"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.486.26"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.978.05"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.17"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.870"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  -20.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.49"
$ws.Range("E24").Value = "  -2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.72"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.69"
$ws.Range("E26").Value = "  -4.95%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.01"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.91"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.08"
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.96"
$ws.Range("E32").Value = "  -5.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.80"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "628.02"
$ws.Range("E35").Value = "  +9.95%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.56"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0989"
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.67"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0471"
$ws.Range("E39").Value = "  +6.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.30"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0733"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.336.02"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.309"
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.94"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.53"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.129"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.05"
$ws.Range("E50").Value = "  -0.40%  "
